$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 629, shifting existing rows 629:712 down to 630:713.
$ws.Rows("629").Insert()

# Populate the newly inserted row 629 with a new price record (same
# market/product metadata as the rest of the table, new date/price/origin).
$ws.Cells.Item(629, 1).Value2 = 9
$ws.Cells.Item(629, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(629, 3).Value2 = "Metropolitana"
$ws.Cells.Item(629, 4).Value2 = 45131
$ws.Cells.Item(629, 5).Value2 = 13
$ws.Cells.Item(629, 6).Value2 = "Fruta"
$ws.Cells.Item(629, 7).Value2 = 100108
$ws.Cells.Item(629, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(629, 9).Value2 = 100108002
$ws.Cells.Item(629, 10).Value2 = "Mango"
$ws.Cells.Item(629, 11).Value2 = "Sin especificar"
$ws.Cells.Item(629, 12).Value2 = "Primera"
$ws.Cells.Item(629, 13).Value2 = 470
$ws.Cells.Item(629, 14).Value2 = 7500
$ws.Cells.Item(629, 15).Value2 = 8000
$ws.Cells.Item(629, 16).Value2 = 7702
$ws.Cells.Item(629, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(629, 18).Value2 = "Brasil"
$ws.Cells.Item(629, 19).Value2 = 1926
$ws.Cells.Item(629, 20).Value2 = 4

# Give the new date cell the same date/time number format used by the
# rest of column D (style index 2 in the original workbook).
$ws.Cells.Item(629, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
